# Auto - Update data with bot!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: update title text
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 32: update title text and link
$ws.Range("D32").Value = "Dynamic Time Warping(DTW)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/352"

# Row 46: update title text and link
$ws.Range("D46").Value = "실신(Syncope) 감별진단"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/447"
